$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 2: the newest meter reading (Number=2, 7/26/2024) ---------------
# Columns: A=Number, B=Id, C=Serial Number, D=Reading Time, E=End Index, F=Voltage, G=Current
$ws.Range("A2:D2").NumberFormat = "@"   # keep these as text so they are not
                                         # auto-detected as numbers/dates

$ws.Range("A2").Value2 = "2"
$ws.Range("B2").Value2 = "d4dfcd97-bc36-4a6f-8425-863a23ecca0b"
$ws.Range("C2").Value2 = "ABC12345"
$ws.Range("D2").Value2 = "7/26/2024"

$ws.Range("A2:D2").Style = "Normal"     # restore default styling/number format

$ws.Range("E2").Value2 = 4
$ws.Range("F2").Value2 = -3
$ws.Range("G2").Value2 = 4

# --- New row 3: the previous reading, pushed down (Number=1, 7/6/2024) ------
$ws.Range("A3:D3").NumberFormat = "@"

$ws.Range("A3").Value2 = "1"
$ws.Range("B3").Value2 = "b0a69cbb-8b0c-41f3-b6b9-fd73c35ee3df"
$ws.Range("C3").Value2 = "ABC12345"
$ws.Range("D3").Value2 = "7/6/2024"

$ws.Range("A3:D3").Style = "Normal"

$ws.Range("E3").Value2 = 0.5
$ws.Range("F3").Value2 = 0.4
$ws.Range("G3").Value2 = 2.3
